# Insert a new data row before current row 41 (shifting existing rows 41-171
# down to 42-172) and populate it with the new weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(41).Insert()

$ws.Range("A41").Value = 7
$ws.Range("B41").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C41").Value = "Ñuble"
$ws.Range("D41").Value = 44459
$ws.Range("E41").Value = 16
$ws.Range("F41").Value = 100114013
$ws.Range("G41").Value = "Zanahoria"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 160
$ws.Range("K41").Value = 8500
$ws.Range("L41").Value = 9000
$ws.Range("M41").Value = 8750
$ws.Range("N41").Value = "$/saco 20 kilos"
$ws.Range("O41").Value = "Provincia de Diguillín"
$ws.Range("P41").Value = 438
$ws.Range("Q41").Value = 20
$ws.Range("R41").Value = "Hortaliza"
